$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

function Set-LiteralText($addr, $text, $styleDonorAddr) {
    # Excel auto-detects date-like strings (e.g. "31-JAN-26") and converts
    # them to a date serial + date-numbered style on plain assignment.
    # Writing with a leading apostrophe keeps it literal text, but stamps
    # the cell with a quotePrefix style variant. Copying an *untouched*
    # donor cell's format (captured before the edit, since assigning a
    # value cancels any pending Copy clipboard) and pasting it back over
    # the edited cell drops that stray style so the cell keeps its
    # original style index.
    $ws.Range($styleDonorAddr).Copy() | Out-Null
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
}

function Set-RowStyleFrom($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# ---- Row 2 ----
Set-RowStyleFrom "A2" "A2"
Set-LiteralText "A2" "31-JAN-26"
$ws.Range("C2").Value = "Nile Air NP-141"
$ws.Range("D2").Value = 9307
$ws.Range("E2").Value = 10072
$ws.Range("F2").Value = -765

# ---- Row 3 ----
Set-LiteralText "A3" "26-MAR-26"
$ws.Range("C3").Value = "Nile Air NP-141"
$ws.Range("D3").Value = 14713
$ws.Range("E3").Value = 14738
$ws.Range("F3").Value = -25

# ---- Row 4 ----
Set-LiteralText "A4" "28-MAR-26"
# C4 unchanged ("Nile Air NP-141")
$ws.Range("D4").Value = 16331
$ws.Range("E4").Value = 16507
$ws.Range("F4").Value = -176

# ---- Row 5 ----
Set-LiteralText "A5" "13-MAY-26"
# C5 unchanged ("Air Arabia Egypt E5-585")
$ws.Range("D5").Value = 7662
$ws.Range("E5").Value = 10586
$ws.Range("F5").Value = -2924
Set-RowStyleFrom "J6" "J5"
$ws.Range("J5").Value = "MEDIUM THREAT - MONITOR"

# ---- Row 6 ----
Set-LiteralText "A6" "14-MAY-26"
$ws.Range("C6").Value = "Nile Air NP-141"
$ws.Range("D6").Value = 8027
$ws.Range("E6").Value = 10586
$ws.Range("F6").Value = -2559
# J6 unchanged ("MEDIUM THREAT - MONITOR")

# ---- Row 7 ----
Set-LiteralText "A7" "16-MAY-26"
$ws.Range("C7").Value = "Nile Air NP-141"
$ws.Range("D7").Value = 8027
$ws.Range("E7").Value = 10586
$ws.Range("F7").Value = -2559
# J7 unchanged ("MEDIUM THREAT - MONITOR")

# ---- New rows 8-11: seed formatting from row 7 (already MEDIUM THREAT styled) ----
Set-RowStyleFrom "A7:K7" "A8:K8"
Set-RowStyleFrom "A7:K7" "A9:K9"
Set-RowStyleFrom "A7:K7" "A10:K10"
Set-RowStyleFrom "A7:K7" "A11:K11"

# ---- Row 8 ----
Set-LiteralText "A8" "20-MAY-26"
$ws.Range("B8").Value = "SM-329"
$ws.Range("C8").Value = "Air Arabia Egypt E5-585"
$ws.Range("D8").Value = 7662
$ws.Range("E8").Value = 10586
$ws.Range("F8").Value = -2924
$ws.Range("G8").Value = 30
$ws.Range("H8").Value = 30
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = "MEDIUM THREAT - MONITOR"
$ws.Range("K8").Value = "EGP"

# ---- Row 9 ----
Set-LiteralText "A9" "21-MAY-26"
$ws.Range("B9").Value = "SM-329"
$ws.Range("C9").Value = "Nile Air NP-141"
$ws.Range("D9").Value = 8027
$ws.Range("E9").Value = 10586
$ws.Range("F9").Value = -2559
$ws.Range("G9").Value = 30
$ws.Range("H9").Value = 30
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "MEDIUM THREAT - MONITOR"
$ws.Range("K9").Value = "EGP"

# ---- Row 10 ----
Set-LiteralText "A10" "28-MAY-26"
$ws.Range("B10").Value = "SM-329"
$ws.Range("C10").Value = "Nile Air NP-141"
$ws.Range("D10").Value = 11013
$ws.Range("E10").Value = 14939
$ws.Range("F10").Value = -3926
$ws.Range("G10").Value = 30
$ws.Range("H10").Value = 30
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = "MEDIUM THREAT - MONITOR"
$ws.Range("K10").Value = "EGP"

# ---- Row 11 ----
Set-LiteralText "A11" "30-MAY-26"
$ws.Range("B11").Value = "SM-329"
$ws.Range("C11").Value = "Nile Air NP-141"
$ws.Range("D11").Value = 11013
$ws.Range("E11").Value = 14939
$ws.Range("F11").Value = -3926
$ws.Range("G11").Value = 30
$ws.Range("H11").Value = 30
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "MEDIUM THREAT - MONITOR"
$ws.Range("K11").Value = "EGP"
